$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "58.762.92"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +1.33%  "

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.001.35"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +2.92%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "563.08"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.84%  "

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "138.52"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +12.18%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.518"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +4.83%  "

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.996.24"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +2.89%  "

# Row 10
$ws.Range("E10").Value = "  +4.64%  "

# Row 11
$ws.Range("E11").Value = "  +2.49%  "

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.453"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +3.94%  "

# Row 13
$ws.Range("E13").Value = "  +7.15%  "

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "33.83"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +5.96%  "

# Row 15
$ws.Range("E15").Value = "  +3.06%  "

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.491.62"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +3.41%  "

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "7.04"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +7.54%  "

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.996.74"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +3.16%  "

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "58.694.42"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.45%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "428.40"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +4.28%  "

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "13.63"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +5.18%  "

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.716"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +8.22%  "

# Row 23
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "13.54"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +6.16%  "

# Row 24
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "7.14"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +4.66%  "

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "80.46"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +3.93%  "

# Row 26
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.29%  "

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.11"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +9.43%  "

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.54"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +2.75%  "

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "7.72"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +5.83%  "

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "25.81"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +3.87%  "

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "6.11"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.35%  "

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.0985"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +2.67%  "

# Row 34
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.970"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +6.12%  "

# Row 35
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "5.78"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +7.46%  "

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0748"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +18.40%  "

# Row 37
$ws.Range("E37").Value = "  +3.56%  "

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "49.02"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +2.49%  "

# Row 39
$ws.Range("E39").Value = "  +4.25%  "

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.78"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +15.16%  "

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "397.07"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +9.90%  "

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.0351"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +1.35%  "

# Row 43
$ws.Range("E43").Value = "  +1.65%  "

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.734.96"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +4.29%  "

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.248"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +7.80%  "

# Row 46
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.06%  "

# Row 47
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "125.47"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +6.42%  "

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.04"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +4.39%  "

# Row 49
$ws.Range("E49").Value = "  +2.60%  "

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "23.53"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +3.41%  "

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "31.42"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +12.50%  "
